$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C4").Value = 0.4786078607860785
$ws.Range("D4").Value = 0.9807140714071406
$ws.Range("E4").Value = -0.9975997599759974
$ws.Range("F4").Value = 0.9814101410141013
$ws.Range("H4").Value = 0.698013801380138
$ws.Range("I4").Value = -0.03313531353135313
$ws.Range("J4").Value = -0.3424422442244224
$ws.Range("K4").Value = 0.1583798379837984
$ws.Range("L4").Value = -0.0722112211221122
$ws.Range("M4").Value = 0.0326072607260726
$ws.Range("N4").Value = 0.3872187218721871
$ws.Range("P4").Value = 0.03977197719771978

$ws.Range("C5").Value = 0.00528052805280528
$ws.Range("D5").Value = -0.1148874887488749
$ws.Range("E5").Value = 0.1164476447644764
$ws.Range("F5").Value = -0.1027542754275428
$ws.Range("H5").Value = -0.2465886588658865
$ws.Range("I5").Value = -0.1426222622262226
$ws.Range("J5").Value = -0.06929492949294928
$ws.Range("K5").Value = 0.1991839183918392
$ws.Range("L5").Value = -0.02035403540354035
$ws.Range("M5").Value = 0.042004200420042
$ws.Range("N5").Value = 0.118967896789679
$ws.Range("P5").Value = 0.01635763576357635

$ws.Range("C6").Value = 0.1398499849984998
$ws.Range("D6").Value = 0.09605760576057605
$ws.Range("E6").Value = -0.1193759375937594
$ws.Range("F6").Value = 0.1609120912091209
$ws.Range("H6").Value = 0.08194419441944194
$ws.Range("I6").Value = -0.04927692769276927
$ws.Range("J6").Value = 0.03259525952595259
$ws.Range("K6").Value = -0.02126612661266126
$ws.Range("L6").Value = -0.1233483348334833
$ws.Range("M6").Value = 0.01782178217821782
$ws.Range("N6").Value = -0.03426342634263426
$ws.Range("P6").Value = 0.01292529252925292

$ws.Range("C7").Value = 0.1547434743474347
$ws.Range("D7").Value = 0.1198679867986799
$ws.Range("E7").Value = 0.04458445844584458
$ws.Range("F7").Value = -0.1604320432043204
$ws.Range("H7").Value = -0.08806480648064806
$ws.Range("I7").Value = 0.007704770477047704
$ws.Range("J7").Value = 0.1778697869786978
$ws.Range("K7").Value = 0.1843864386438644
$ws.Range("L7").Value = 0.9638283828382838
$ws.Range("M7").Value = -0.06255025502550254
$ws.Range("N7").Value = -0.4957815781578156
$ws.Range("P7").Value = -0.06159015901590158

$ws.Range("C8").Value = 0.04123612361236123
$ws.Range("D8").Value = -0.06366636663666365
$ws.Range("E8").Value = 0.07192319231923192
$ws.Range("F8").Value = -0.01816981698169817
$ws.Range("H8").Value = 0.6435643564356436
$ws.Range("I8").Value = 0.8634503450345034
$ws.Range("J8").Value = -0.1334293429342934
$ws.Range("K8").Value = -0.00096009600960096
$ws.Range("L8").Value = 0.004872487248724872
$ws.Range("M8").Value = -0.1418781878187819
$ws.Range("N8").Value = -0.0257065706570657
$ws.Range("P8").Value = 0.1338493849384938

$ws.Range("C9").Value = 0.4663306330633062
$ws.Range("D9").Value = 0.02339033903390339
$ws.Range("E9").Value = -0.03906390639063906
$ws.Range("F9").Value = 0.02047404740474047
$ws.Range("H9").Value = 0.02973897389738973
$ws.Range("I9").Value = 0.09438943894389439
$ws.Range("J9").Value = -0.1061386138613861
$ws.Range("K9").Value = -0.6338433843384337
$ws.Range("L9").Value = -0.03203120312031203
$ws.Range("M9").Value = -0.1226162616261626
$ws.Range("N9").Value = 0.1318811881188119
$ws.Range("P9").Value = 0.06294629462946294

$ws.Range("C10").Value = 0.4871767176717672
$ws.Range("D10").Value = -0.04802880288028803
$ws.Range("E10").Value = 0.06327032703270327
$ws.Range("F10").Value = -0.0571137113711371
$ws.Range("H10").Value = -0.0114971497149715
$ws.Range("I10").Value = -0.03294329432943294
$ws.Range("J10").Value = -0.006948694869486947
$ws.Range("K10").Value = 0.6153975397539753
$ws.Range("L10").Value = 0.07278727872787277
$ws.Range("M10").Value = -0.07955595559555954
$ws.Range("N10").Value = -0.05945394539453944
$ws.Range("P10").Value = -0.04014401440144014

$ws.Range("C11").Value = -0.1402820282028203
$ws.Range("D11").Value = 0.0919051905190519
$ws.Range("E11").Value = -0.1138553855385538
$ws.Range("F11").Value = 0.1194959495949595
$ws.Range("H11").Value = -0.06720672067206719
$ws.Range("I11").Value = -0.1291929192919292
$ws.Range("J11").Value = -0.1872787278727873
$ws.Range("K11").Value = -0.03692769276927692
$ws.Range("L11").Value = -0.1444944494449445
$ws.Range("M11").Value = -0.06367836783678367
$ws.Range("N11").Value = 0.1699129912991299
$ws.Range("P11").Value = -0.04688868886888688

$ws.Range("C12").Value = -0.02157815781578157
$ws.Range("D12").Value = -0.1859345934593459
$ws.Range("E12").Value = 0.1747254725472547
$ws.Range("F12").Value = -0.1805460546054605
$ws.Range("H12").Value = -0.1992559255925593
$ws.Range("I12").Value = -0.1041584158415841
$ws.Range("J12").Value = 0.1633843384338434
$ws.Range("K12").Value = 0.04046804680468046
$ws.Range("L12").Value = 0.04136813681368137
$ws.Range("M12").Value = -0.1804860486048605
$ws.Range("N12").Value = -0.08296429642964295
$ws.Range("P12").Value = -0.07615961596159614

$ws.Range("C13").Value = 0.1688688868886888
$ws.Range("D13").Value = 0.07761176117611761
$ws.Range("E13").Value = -0.08174017401740175
$ws.Range("F13").Value = 0.05948994899489948
$ws.Range("H13").Value = 0.03918391839183918
$ws.Range("I13").Value = 0.4652265226522652
$ws.Range("J13").Value = -0.01045304530453045
$ws.Range("K13").Value = -0.2018361836183618
$ws.Range("L13").Value = -0.004908490849084907
$ws.Range("M13").Value = 0.08784878487848785
$ws.Range("N13").Value = 0.07075907590759076
$ws.Range("P13").Value = 0.09366936693669366

$ws.Range("C14").Value = -0.2346234623462346
$ws.Range("D14").Value = -0.01711371137113711
$ws.Range("E14").Value = -0.01274527452745275
$ws.Range("F14").Value = 0.004584458445844585
$ws.Range("H14").Value = 0.04301230123012301
$ws.Range("I14").Value = -0.03866786678667866
$ws.Range("J14").Value = -0.1127272727272727
$ws.Range("K14").Value = -0.09989798979897989
$ws.Range("L14").Value = -0.1572397239723972
$ws.Range("M14").Value = -0.03512751275127513
$ws.Range("N14").Value = 0.115031503150315
$ws.Range("P14").Value = -0.1078667866786679

